$d = $word.ActiveDocument

# --- Hunk 1: "ENC/3 PRACOWNIK (...)" paragraph ---
# wrap IdPrac in proofErr spellStart/spellEnd, split the trailing run so
# " Stanowisko," is its own run between "Imie," and " Wynagrodzenie)"
$p1 = $d.Paragraphs.Item(5)
if ($p1.Range.Text -notmatch "ENC/3 PRACOWNIK") {
    throw "Paragraph 5 is not the expected 'ENC/3 PRACOWNIK' paragraph: $($p1.Range.Text)"
}
$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00093804" w:rsidRDefault="00093804" w:rsidP="00093804"><w:pPr><w:ind w:left="708"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">ENC/3 PRACOWNIK </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00AC1FD3"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr><w:t>IdPrac</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>, Hasło, Nazwisko, Imię,</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Stanowisko,</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Wynagrodzenie)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p1.Range.InsertXML($xml1)

# --- Hunk 2: "Pracownicy(IdPrac, ...)" paragraph ---
# split the trailing run the same way, and move the _GoBack bookmark here
# (between " Stanowisko," and " Wynagrodzenie)")
$p2 = $d.Paragraphs.Item(8)
if ($p2.Range.Text -notmatch "Pracownicy\(") {
    throw "Paragraph 8 is not the expected 'Pracownicy(' paragraph: $($p2.Range.Text)"
}
$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00093804" w:rsidRDefault="00093804" w:rsidP="00093804"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/><w:t>Pracownicy(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00093804"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr><w:t>IdPrac</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>, Hasło, Nazwisko, Imię,</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Stanowisko,</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Wynagrodzenie)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p2.Range.InsertXML($xml2)

# --- Hunk 3: "Klienci (IdKlienta, ...)" paragraph ---
# merge the tab / "Klienci" / " " / "(" runs into a single tab run whose
# text is "Klienci ("
$p3 = $d.Paragraphs.Item(33)
if ($p3.Range.Text -notmatch "Klienci") {
    throw "Paragraph 33 is not the expected 'Klienci (' paragraph: $($p3.Range.Text)"
}
$xml3 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00D611F7" w:rsidRDefault="00D611F7" w:rsidP="00D611F7"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/><w:t>Klienci (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr><w:t>IdKlienta</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>, Nazwisko, Imię, Telefon, Email)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p3.Range.InsertXML($xml3)

# --- Hunk 4: "Rachunki(IdRachunku, ..., Suma)" paragraph ---
# remove the _GoBack bookmark that used to sit at the end of this paragraph
# (it moved up to the "Pracownicy(" paragraph, hunk 2)
$p4 = $d.Paragraphs.Item(41)
if ($p4.Range.Text -notmatch "Rachunki\(") {
    throw "Paragraph 41 is not the expected 'Rachunki(' paragraph: $($p4.Range.Text)"
}
$xml4 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00D611F7" w:rsidRDefault="009F6E4D" w:rsidP="00093804"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/><w:t>Rachunki(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr><w:t>IdRachunku</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>DataWystawienia</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>DataOpłacenia</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>, Suma)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p4.Range.InsertXML($xml4)

Write-Output "done"
